$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.699.75'
$ws.Range('E2').Value = '  -2.52%  '
$ws.Range('D3').Value = '1.557.01'
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('E4').Value = '  +0.05%  '
$c = $ws.Range('D5')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '205.69'
$c.Style = $origStyle
$ws.Range('E5').Value = '  -1.22%  '
$c = $ws.Range('D6')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.489'
$c.Style = $origStyle
$ws.Range('E6').Value = '  -1.91%  '
$c = $ws.Range('D8')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '21.96'
$c.Style = $origStyle
$ws.Range('E8').Value = '  +0.32%  '
$ws.Range('E9').Value = '  -0.63%  '
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('D12').Value = '1.778.04'
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('D13').Value = '1.564.05'
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('E14').Value = '  -2.32%  '
$ws.Range('E15').Value = '  -0.82%  '
$c = $ws.Range('D16')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '61.57'
$c.Style = $origStyle
$ws.Range('E16').Value = '  -2.80%  '
$ws.Range('D17').Value = '26.736.05'
$ws.Range('E17').Value = '  -2.32%  '
$c = $ws.Range('D18')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '7.36'
$c.Style = $origStyle
$ws.Range('E18').Value = '  +1.35%  '
$c = $ws.Range('D19')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '213.25'
$c.Style = $origStyle
$ws.Range('E19').Value = '  +0.44%  '
$ws.Range('E20').Value = '  -1.95%  '
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('E22').Value = '  -0.53%  '
$c = $ws.Range('D23')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '9.34'
$c.Style = $origStyle
$ws.Range('E23').Value = '  -1.85%  '
$c = $ws.Range('D24')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.00'
$c.Style = $origStyle
$ws.Range('E24').Value = '  +0.15%  '
$c = $ws.Range('D25')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '152.39'
$c.Style = $origStyle
$ws.Range('E25').Value = '  -0.67%  '
$c = $ws.Range('D26')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '6.78'
$c.Style = $origStyle
$ws.Range('E26').Value = '  +0.83%  '
$c = $ws.Range('D27')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '14.81'
$c.Style = $origStyle
$ws.Range('E27').Value = '  -1.23%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('E30').Value = '  -1.53%  '
$ws.Range('E31').Value = '  -3.87%  '
$c = $ws.Range('D32')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.14'
$c.Style = $origStyle
$ws.Range('E32').Value = '  -1.74%  '
$ws.Range('D33').Value = '1.385.07'
$ws.Range('E33').Value = '  +1.59%  '
$ws.Range('E34').Value = '  -1.34%  '
$ws.Range('E35').Value = '  +1.00%  '
$ws.Range('E36').Value = '  -0.79%  '
$c = $ws.Range('D37')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.931'
$c.Style = $origStyle
$ws.Range('E37').Value = '  -4.28%  '
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('E39').Value = '  -2.61%  '
$c = $ws.Range('D40')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.812'
$c.Style = $origStyle
$ws.Range('E40').Value = '  -1.29%  '
$ws.Range('E41').Value = '  +0.06%  '
$c = $ws.Range('D42')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.993'
$c.Style = $origStyle
$ws.Range('E42').Value = '  +2.11%  '
$ws.Range('E43').Value = '  +1.95%  '
$c = $ws.Range('D44')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.17'
$c.Style = $origStyle
$ws.Range('E44').Value = '  +1.18%  '
$ws.Range('E45').Value = '  -1.67%  '
$c = $ws.Range('D46')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '63.06'
$c.Style = $origStyle
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('D47').Value = '1.691.02'
$ws.Range('E47').Value = '  -0.50%  '
$c = $ws.Range('D48')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '85.46'
$c.Style = $origStyle
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('D49').Value = '0.0₇0972'
$ws.Range('E49').Value = '  -2.15%  '
$c = $ws.Range('D50')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0493'
$c.Style = $origStyle
$ws.Range('E50').Value = '  -0.23%  '
$c = $ws.Range('D51')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0947'
$c.Style = $origStyle
$ws.Range('E51').Value = '  -0.88%  '
